# Adapt column header formatting to respective input file names (#7)
# - rename the "_old"/"_new" header suffixes to "_FV2404"/"_FV2410"
# - turn the header range into a proper Excel Table (ListObject)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header row (row 1, columns A:U) from the old "_old"/"_new"
#    suffix convention to the new "_FV2404"/"_FV2410" format-version suffix.
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Wrap the used range in a native Excel Table so the header row gets
#    table/autofilter semantics (mirrors the exporter's `Table1` over A1:U79).
$usedRange = $ws.Range("A1:U79")
$lo = $ws.ListObjects.Add(1, $usedRange, $null, 1)
$lo.Name = "Table1"

# 3. Freeze the header row so it stays visible while scrolling.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
